# Daily attendance processing - 2025-11-09 09:20:52
# Swap the order of "System" and the other recorder (email) in column G
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "System, *") {
        $other = $val.Substring(8)
        $cell.Value2 = "$other, System"
    }
}
